# Add the new "Sheet 2" worksheet right after "Worksheet"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Worksheet")

$sheet2 = $wb.Worksheets.Add($null, $ws)
$sheet2.Name = "Sheet 2"

# Populate the new summary table (mahasiswa asing / TS.. breakdown)
$sheet2.Range("A1").Value = "TS"
$sheet2.Range("C1").Value = 801
$sheet2.Range("D1").Value = 52

$sheet2.Range("A2").Value = "TS-1"
$sheet2.Range("C2").Value = 666
$sheet2.Range("G2").Value = 868

$sheet2.Range("A3").Value = "TS-2"
$sheet2.Range("B3").Value = 232
$sheet2.Range("C3").Value = 4328
$sheet2.Range("D3").Value = 232
$sheet2.Range("E3").Value = 192
$sheet2.Range("G3").Value = 838

$sheet2.Range("A4").Value = "TS-3"
$sheet2.Range("B4").Value = 227
$sheet2.Range("C4").Value = 4060
$sheet2.Range("D4").Value = 227
$sheet2.Range("E4").Value = 203
$sheet2.Range("G4").Value = 911

$sheet2.Range("A5").Value = "TS-4"
$sheet2.Range("B5").Value = 297
$sheet2.Range("C5").Value = 4937
$sheet2.Range("D5").Value = 297
$sheet2.Range("E5").Value = 270
$sheet2.Range("G5").Value = 959

# Apply the additional autofilter on the "Tahun Akademik" column (field 4 of
# the B1:L17 filter range == column E) restricting to TS / TS-1 .. TS-4.
# This also hides row 17 (TS-5), matching the recorded change.
$ws.Range("B1:L17").AutoFilter(4, @("TS-4", "TS-3", "TS-2", "TS-1", "TS"), 7) | Out-Null

# Make sure row 17 is explicitly marked hidden as well.
$ws.Rows.Item(17).Hidden = $true

# Keep the original sheet ("Worksheet") as the active/selected tab.
$ws.Activate()
